# Adding in the raw data as downloaded via the RPP script
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update J2's note text (dropped the "(verify)" suffix)
$ws.Range("J2").Value = "Not significant original"

# Insert a new column before column K. This shifts the existing
# K:P header cells (original_r_collected .. replication_n_collected)
# one column to the right, to L:Q, making room for a new
# "Data_Cleaned" column.
$ws.Columns("K:K").Insert()

# New column header
$ws.Range("K1").Value = "Data_Cleaned"

# Mark rows 2, 3 and 5 as having cleaned data
$ws.Range("K2").Value = 1
$ws.Range("K3").Value = 1
$ws.Range("K5").Value = 1

# Row 4: remove the stray n_studies value that didn't belong
$ws.Range("H4").ClearContents()

# Row 5: the "1" flag that used to live in column I moved into the
# new Data_Cleaned column (K5, set above) - clear the old cell
$ws.Range("I5").ClearContents()

# Row 8: remove the old "6 / 6 with effect sizes not convertible to r" note
$ws.Range("I8").ClearContents()
$ws.Range("J8").ClearContents()

# Put the selection where it ended up in the authored workbook
$ws.Range("J6").Select()
